$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date values from 45221 (2023-10-22) to 45224 (2023-10-25)
# for rows 2-5, preserving existing cell formatting.
$ws.Range("C2").Value = 45224
$ws.Range("C3").Value = 45224
$ws.Range("C4").Value = 45224
$ws.Range("C5").Value = 45224
